# Refresh cryptos list: update Price (D) and Volume(1h) (E) figures,
# and fix two rows where the coin order was swapped back
# (Chainlink/Uniswap at 19-20, BabyDogeCoin/ARBITRUM at 47-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.599.80"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "2.653.73"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'597.34"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").Value = "'168.28"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.545"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").Value = "2.654.76"
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "'5.27"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").Value = "'28.10"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "3.137.85"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("E16").Value = "  -3.30%  "
$ws.Range("D17").Value = "67.603.18"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "2.665.53"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'12.08"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'8.27"
$ws.Range("E20").Value = "  +7.83%  "
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").Value = "'4.80"
$ws.Range("E23").Value = "  -3.95%  "
$ws.Range("D24").Value = "'11.03"
$ws.Range("E24").Value = "  +8.55%  "
$ws.Range("D25").Value = "'2.01"
$ws.Range("E25").Value = "  -4.01%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'71.10"
$ws.Range("E27").Value = "  -3.46%  "
$ws.Range("D28").Value = "2.793.43"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").Value = "'560.35"
$ws.Range("E31").Value = "  -5.30%  "
$ws.Range("E32").Value = "  -2.91%  "
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("D34").Value = "'1.93"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -4.51%  "
$ws.Range("D38").Value = "'158.15"
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").Value = "'19.42"
$ws.Range("E39").Value = "  -2.95%  "
$ws.Range("D40").Value = "'0.373"
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("D42").Value = "'1.83"
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("E44").Value = "  -4.62%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'40.25"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0301"
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.598"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").Value = "'154.66"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").Value = "'3.89"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("E51").Value = "  -2.84%  "
